$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J11").Value = 44460
$ws.Range("J11").NumberFormat = $ws.Range("J7").NumberFormat
$ws.Range("K11").Value = 1746
$ws.Range("L11").Value = "36/24"
$ws.Range("M11").Value = 82
$ws.Range("N11").Value = 645

$ws.Range("N12").Select()
